# Roteiro de entregas - caminhão 1
# Replace the delivery schedule data (rows 2-18 originally) with the
# updated route: dates shift forward (19/11 -> 21/11, 21/11 -> 22/11),
# order numbers/neighbourhoods are reassigned, and six new rows (19-24)
# are appended for the extra 22/11 afternoon deliveries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('11197','21/11/2024','manhã','Desconhecido','Picadas do Sul'),
    @('11157','21/11/2024','manhã','Desconhecido','Forquilhinha'),
    @('11071','21/11/2024','manhã','Desconhecido','Barreiros'),
    @('11073','21/11/2024','manhã','Desconhecido','Bela Vista'),
    @('11165','21/11/2024','manhã','Desconhecido','Ipiranga'),
    @('11167','21/11/2024','manhã','Desconhecido','Ipiranga'),
    @('11200','21/11/2024','tarde','Desconhecido','Ponta de Baixo'),
    @('11201','21/11/2024','tarde','Desconhecido','Ponte do Imaruim'),
    @('11144','21/11/2024','tarde','Desconhecido','Centro'),
    @('11143','21/11/2024','tarde','Desconhecido','Centro'),
    @('11142','21/11/2024','tarde','Desconhecido','Centro'),
    @('11181','21/11/2024','tarde','Desconhecido','Pachecos'),
    @('11149','22/11/2024','manhã','Desconhecido','Distrito Industrial'),
    @('11170','22/11/2024','manhã','Desconhecido','Jardim Eldorado'),
    @('11171','22/11/2024','manhã','Desconhecido','Jardim Eldorado'),
    @('11194','22/11/2024','manhã','Desconhecido','Pedra Branca'),
    @('11182','22/11/2024','manhã','Desconhecido','Pagani'),
    @('11206','22/11/2024','tarde','Desconhecido','Praia Comprida'),
    @('11139','22/11/2024','tarde','Desconhecido','Centro'),
    @('11175','22/11/2024','tarde','Desconhecido','Monte Cristo'),
    @('11116','22/11/2024','tarde','Desconhecido','Campinas'),
    @('11177','22/11/2024','tarde','Desconhecido','Nossa Senhora do Rosário'),
    @('11072','22/11/2024','tarde','Desconhecido','Bela Vista')
)

# Template row (row 2) already carries the right style (centered, bordered
# text cells) - extend that style down to the newly added rows (19-24)
# before writing values, so the appended rows look like the rest of the table.
$templateRange = $ws.Range("A2:E2")
$newRowsRange = $ws.Range("A19:E24")
$templateRange.Copy() | Out-Null
$newRowsRange.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Column A holds order numbers stored as text (not numbers) in the source
# data - force the whole column to text format up front (one single style
# change reused by every row) so "11197" etc. keep their original string
# type instead of being auto-converted to a number.
$ws.Range("A2:A24").NumberFormat = "@"

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}
